# "break out stock.yaml completed"
#
# 1. On sheet "10per change": D20:D28 (bsecode) were stored as text and
#    become numeric. A new block of rows (29:37) is appended, duplicating
#    rows 20:28 but stamped with a later "Date Time" (17:42:30 instead of
#    17:18:12) and with bsecode stored back as text, like the original
#    rows were before being "completed" to numeric.
# 2. On sheet "DND 3 V 0.3" the same thing happens for the single row 4,
#    appended as row 5.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "10per change": fix D20:D28 to numeric, append rows 29:37
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("10per change")

$bsecodes1 = @{ 20 = 500510; 21 = 543287; 22 = 542066; 23 = 500112; 24 = 533096; 25 = 500093; 26 = 543396; 27 = 532898; 28 = 532155 }
foreach ($r in $bsecodes1.Keys) {
    $ws1.Range("D$r").Value = $bsecodes1[$r]
}

$rows1 = @(
    @(1, "LT",         "Larsen & Toubro Limited",                 "500510", 0.17,  3409,   10372458),
    @(2, "LODHA",       "Macrotech Developers Ltd",                "543287", -0.53, 1296.85, 1797849),
    @(3, "ATGL",        "Adani Total Gas Ltd",                     "542066", 3.03,  936.25,  6670432),
    @(4, "SBIN",        "State Bank Of India",                     "500112", 1.88,  789.75,  74256082),
    @(5, "ADANIPOWER",  "Adani Power Limited",                     "533096", 0.51,  726.65,  34537620),
    @(6, "CGPOWER",     "CG Power and Industrial Solutions Ltd",   "500093", 0.1,   627.65,  10674892),
    @(7, "PAYTM",       "One 97 Communications Ltd",               "543396", -4.91, 339.85,  7404922),
    @(8, "POWERGRID",   "Power Grid Corporation Of India Limited", "532898", 0.96,  298.8,   45312613),
    @(9, "GAIL",        "Gail (india) Limited",                    "532155", 2.55,  195.15,  49797002)
)

$destRow = 29
foreach ($row in $rows1) {
    $ws1.Cells.Item($destRow, 1).Value = $row[0]
    $ws1.Cells.Item($destRow, 2).Value = $row[1]
    $ws1.Cells.Item($destRow, 3).Value = $row[2]

    $dCell = $ws1.Cells.Item($destRow, 4)
    $dCell.Value = "'" + $row[3]
    $dCell.Style = "Normal"

    $ws1.Cells.Item($destRow, 5).Value = $row[4]
    $ws1.Cells.Item($destRow, 6).Value = $row[5]
    $ws1.Cells.Item($destRow, 7).Value = $row[6]
    $ws1.Cells.Item($destRow, 8).Value = "05/06/2024 17:42:30"

    $destRow++
}

# ---------------------------------------------------------------------
# Sheet "DND 3 V 0.3": fix D4 to numeric, append row 5
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("DND 3 V 0.3")

$ws2.Range("D4").Value = 500331

$d5 = $ws2.Cells.Item(5, 4)
$ws2.Cells.Item(5, 1).Value = 1
$ws2.Cells.Item(5, 2).Value = "PIDILITIND"
$ws2.Cells.Item(5, 3).Value = "Pidilite Industries Limited"
$d5.Value = "'500331"
$d5.Style = "Normal"
$ws2.Cells.Item(5, 5).Value = 3.4
$ws2.Cells.Item(5, 6).Value = 3166.2
$ws2.Cells.Item(5, 7).Value = 632880
$ws2.Cells.Item(5, 8).Value = "05/06/2024 17:42:30"

Write-Output "edit complete"
